$d = $word.ActiveDocument

# --- Paragraph 2: "Hola" -> "Buenos " + proofErr(spellStart) + "dias" + proofErr(spellEnd) ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($ptext -eq "Hola") {
        $target = $d.Paragraphs($i)
        break
    }
}

if ($target -ne $null) {
    $rng = $target.Range

    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">Buenos </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>dias</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

    $rng.InsertXML($xml)
}

# --- Paragraph 3: "bonjour" -> "arigato" ---
$d.Content.Find.Execute("bonjour", $true, $false, $false, $false, $false,
                         $true, 1, $false, "arigato", 2)
